# Added 4wk low sales check - updated forecast figures and derived metrics
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---
$ws1.Range("D2").Value = 331
$ws1.Range("H2").Value = 6.59
$ws1.Range("L2").Value = 1.16

$ws1.Range("D3").Value = 408
$ws1.Range("H3").Value = 4.53
$ws1.Range("L3").Value = 1.17

$ws1.Range("D4").Value = 452
$ws1.Range("H4").Value = 3.19
$ws1.Range("L4").Value = 1.02

$ws1.Range("D5").Value = 412
$ws1.Range("H5").Value = 2.4
$ws1.Range("L5").Value = 1.04

$ws1.Range("D6").Value = 320
$ws1.Range("H6").Value = 1.81
$ws1.Range("L6").Value = 0.89

$ws1.Range("D7").Value = 264
$ws1.Range("H7").Value = 0.98
$ws1.Range("J7").Value = "Urgent"
$ws1.Range("L7").Value = 1.15

$ws1.Range("D8").Value = 301
$ws1.Range("H8").Value = 0
$ws1.Range("I8").Value = "High"
$ws1.Range("J8").Value = "Urgent"
$ws1.Range("L8").Value = 1.04

$ws1.Range("D9").Value = 388
$ws1.Range("H9").Value = 0
$ws1.Range("L9").Value = 0.9399999999999999

$ws1.Range("D10").Value = 428
$ws1.Range("L10").Value = 0.87

$ws1.Range("D11").Value = 382
$ws1.Range("L11").Value = 1.1

$ws1.Range("D12").Value = 314
$ws1.Range("L12").Value = 0.87

$ws1.Range("D13").Value = 310
$ws1.Range("L13").Value = 0.91

$ws1.Range("D14").Value = 380
$ws1.Range("L14").Value = 1.02

$ws1.Range("D15").Value = 439
$ws1.Range("L15").Value = 1.13

$ws1.Range("D16").Value = 418
$ws1.Range("L16").Value = 1.02

$ws1.Range("D17").Value = 349
$ws1.Range("L17").Value = 0.84

# --- Summary sheet ---
# Force text number format so these remain text cells (matching original inlineStr type)
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "5902"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "2879"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "1605"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "452"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "264"
